# Apply the weekly Fruta/Hortaliza update: reshuffle the "Fecha" (D),
# "Calidad" (I), "Volumen" (J), "Precio mínimo" (K), "Precio máximo" (L),
# "Precio promedio ponderado" (M) and "Precio $/Kg" (P) values across
# rows 2-11 as published in the new weekly export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "2022-03-25"

# Row 3
$ws.Range("D3").Value = "2022-03-10"

# Row 4
$ws.Range("D4").Value = "2022-03-17"

# Row 5
$ws.Range("D5").Value = "2022-03-29"

# Row 6
$ws.Range("D6").Value = "2022-03-31"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range("P6").Value = 861

# Row 7
$ws.Range("D7").Value = "2022-03-22"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 100
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 15500
$ws.Range("P7").Value = 861

# Row 8
$ws.Range("D8").Value = "2022-03-08"

# Row 9
$ws.Range("D9").Value = "2022-03-15"

# Row 10
$ws.Range("D10").Value = "2021-07-07"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 17000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 17500
$ws.Range("P10").Value = 972

# Row 11
$ws.Range("D11").Value = "2021-07-07"
$ws.Range("I11").Value = "Segunda"
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 15000
$ws.Range("P11").Value = 833
